# Case_4_177 (380 kV) update: refresh computed pl_mw results in B2:L25
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build the replacement block as a 2-D array (24 rows x 11 cols: B..L)
$data = New-Object "object[,]" 24,11
# Row 2 (A2=0)
$data[0,0] = 3.104166210929861
$data[0,1] = 0.1070510703279695
$data[0,2] = 0.4751578269798671
$data[0,3] = 0.1358606409969916
$data[0,4] = 0
$data[0,5] = 3.718735986197402
$data[0,6] = 2.732116817893342
$data[0,7] = 0
$data[0,8] = 0.04135229587989553
$data[0,9] = 0
$data[0,10] = 0.5765860246627597
# Row 3 (A3=1)
$data[1,0] = 3.021990669697857
$data[1,1] = 0.09358758368014719
$data[1,2] = 0.4740562597430227
$data[1,3] = 0.1362274353605031
$data[1,4] = 0
$data[1,5] = 3.670681800054581
$data[1,6] = 2.717013883749246
$data[1,7] = 0
$data[1,8] = 0.04113641471810769
$data[1,9] = 0
$data[1,10] = 0.5711437283538174
# Row 4 (A4=2)
$data[2,0] = 2.973545632776791
$data[2,1] = 0.08535406081310271
$data[2,2] = 0.4735753984476929
$data[2,3] = 0.1364890244904036
$data[2,4] = 0
$data[2,5] = 3.642943213616405
$data[2,6] = 2.70886909298099
$data[2,7] = 0
$data[2,8] = 0.04099979563715372
$data[2,9] = 0
$data[2,10] = 0.5681073514895445
# Row 5 (A5=3)
$data[3,0] = 2.954309131880734
$data[3,1] = 0.08200701827757939
$data[3,2] = 0.4734285704765426
$data[3,3] = 0.1366047744477612
$data[3,4] = 0
$data[3,5] = 3.632080969258482
$data[3,6] = 2.705832626601421
$data[3,7] = 0
$data[3,8] = 0.04094309488491366
$data[3,9] = 0
$data[3,10] = 0.566946708174541
# Row 6 (A6=4)
$data[4,0] = 2.951145428963571
$data[4,1] = 0.08145173476222567
$data[4,2] = 0.473407155989463
$data[4,3] = 0.1366245473936907
$data[4,4] = 0
$data[4,5] = 3.630303881932917
$data[4,6] = 2.705345458358863
$data[4,7] = 0
$data[4,8] = 0.04093361758331415
$data[4,9] = 0
$data[4,10] = 0.5667586158630797
# Row 7 (A7=5)
$data[5,0] = 2.973284157473529
$data[5,1] = 0.08530888846684093
$data[5,2] = 0.4735732194046136
$data[5,3] = 0.1364905484811878
$data[5,4] = 0
$data[5,5] = 3.642794937928358
$data[5,6] = 2.708826999461081
$data[5,7] = 0
$data[5,8] = 0.040999035116843
$data[5,9] = 0
$data[5,10] = 0.5680913881265184
# Row 8 (A8=6)
$data[6,0] = 3.075414358624698
$data[6,1] = 0.1024018965026698
$data[6,2] = 0.4747374066833316
$data[6,3] = 0.1359795636712864
$data[6,4] = 0
$data[6,5] = 3.701798451779183
$data[6,6] = 2.726674400826937
$data[6,7] = 0
$data[6,8] = 0.04127870134942491
$data[6,9] = 0
$data[6,10] = 0.5746461422154141
# Row 9 (A9=7)
$data[7,0] = 3.291686900319746
$data[7,1] = 0.1361923589979881
$data[7,2] = 0.4785738972673528
$data[7,3] = 0.1352660634907767
$data[7,4] = 0
$data[7,5] = 3.831661117703646
$data[7,6] = 2.770683840174684
$data[7,7] = 0
$data[7,8] = 0.04179508073577587
$data[7,9] = 0
$data[7,10] = 0.5899254255013915
# Row 10 (A10=8)
$data[8,0] = 3.460412948906878
$data[8,1] = 0.161198027986984
$data[8,2] = 0.4823438296671867
$data[8,3] = 0.1349177229060352
$data[8,4] = 0
$data[8,5] = 3.93590580236048
$data[8,6] = 2.808595873223481
$data[8,7] = 0
$data[8,8] = 0.04215526569708317
$data[8,9] = 0
$data[8,10] = 0.6026373009228507
# Row 11 (A11=9)
$data[9,0] = 3.539325108957314
$data[9,1] = 0.1726162168465351
$data[9,2] = 0.4842664409970467
$data[9,3] = 0.1347974457481662
$data[9,4] = 0
$data[9,5] = 3.985292964906876
$data[9,6] = 2.827073271420431
$data[9,7] = 0
$data[9,8] = 0.04231502900741546
$data[9,9] = 0
$data[9,10] = 0.6087448649761313
# Row 12 (A12=10)
$data[10,0] = 3.569518609519662
$data[10,1] = 0.1769464013268589
$data[10,2] = 0.4850244132779125
$data[10,3] = 0.1347573906856407
$data[10,4] = 0
$data[10,5] = 4.004280716410676
$data[10,6] = 2.83424864685054
$data[10,7] = 0
$data[10,8] = 0.04237494517206564
$data[10,9] = 0
$data[10,10] = 0.611104477439369
# Row 13 (A13=11)
$data[11,0] = 3.563002044412258
$data[11,1] = 0.1760135325757801
$data[11,2] = 0.4848598386302996
$data[11,3] = 0.134765773024581
$data[11,4] = 0
$data[11,5] = 4.000178596798605
$data[11,6] = 2.832695345659943
$data[11,7] = 0
$data[11,8] = 0.0423620669920961
$data[11,9] = 0
$data[11,10] = 0.6105942091835743
# Row 14 (A14=12)
$data[12,0] = 3.541802904646261
$data[12,1] = 0.172972335104788
$data[12,2] = 0.4843281998071802
$data[12,3] = 0.1347940403372121
$data[12,4] = 0
$data[12,5] = 3.986849350036152
$data[12,6] = 2.827660010874524
$data[12,7] = 0
$data[12,8] = 0.0423199699967185
$data[12,9] = 0
$data[12,10] = 0.60893805299402
# Row 15 (A15=13)
$data[13,0] = 3.528858371386093
$data[13,1] = 0.1711103482123235
$data[13,2] = 0.4840064541557467
$data[13,3] = 0.1348120700263102
$data[13,4] = 0
$data[13,5] = 3.978722131386007
$data[13,6] = 2.824598994854512
$data[13,7] = 0
$data[13,8] = 0.04229410862968308
$data[13,9] = 0
$data[13,10] = 0.6079297085489799
# Row 16 (A16=14)
$data[14,0] = 3.455299348876395
$data[14,1] = 0.1604527029185476
$data[14,2] = 0.4822223661185205
$data[14,3] = 0.1349263517003685
$data[14,4] = 0
$data[14,5] = 3.932718066438753
$data[14,6] = 2.807413219510181
$data[14,7] = 0
$data[14,8] = 0.04214474297418214
$data[14,9] = 0
$data[14,10] = 0.6022447023804745
# Row 17 (A17=15)
$data[15,0] = 3.410726552409528
$data[15,1] = 0.1539257198647874
$data[15,2] = 0.4811811127629966
$data[15,3] = 0.1350062399986331
$data[15,4] = 0
$data[15,5] = 3.905001824636003
$data[15,6] = 2.797186564101366
$data[15,7] = 0
$data[15,8] = 0.04205206877331769
$data[15,9] = 0
$data[15,10] = 0.5988404071107283
# Row 18 (A18=16)
$data[16,0] = 3.385292507421013
$data[16,1] = 0.1501756063244102
$data[16,2] = 0.4806017525541506
$data[16,3] = 0.1350557837221285
$data[16,4] = 0
$data[16,5] = 3.889245004014214
$data[16,6] = 2.791420249519575
$data[16,7] = 0
$data[16,8] = 0.0419983800713446
$data[16,9] = 0
$data[16,10] = 0.596912918556086
# Row 19 (A19=17)
$data[17,0] = 3.376715824723192
$data[17,1] = 0.1489065697950309
$data[17,2] = 0.4804089455388407
$data[17,3] = 0.1350731756567871
$data[17,4] = 0
$data[17,5] = 3.88394166176181
$data[17,6] = 2.78948771967481
$data[17,7] = 0
$data[17,8] = 0.0419801357356242
$data[17,9] = 0
$data[17,10] = 0.5962655520469866
# Row 20 (A20=18)
$data[18,0] = 3.41545038173291
$data[18,1] = 0.1546201094986372
$data[18,2] = 0.4812899330873393
$data[18,3] = 0.1349973637770088
$data[18,4] = 0
$data[18,5] = 3.907933113516634
$data[18,6] = 2.798263215666196
$data[18,7] = 0
$data[18,8] = 0.04206197390575817
$data[18,9] = 0
$data[18,10] = 0.5991996353983637
# Row 21 (A21=19)
$data[19,0] = 3.548021153205298
$data[19,1] = 0.1738654343031101
$data[19,2] = 0.4844835424646448
$data[19,3] = 0.1347855885059648
$data[19,4] = 0
$data[19,5] = 3.990756687581779
$data[19,6] = 2.829134158443935
$data[19,7] = 0
$data[19,8] = 0.04233235068027241
$data[19,9] = 0
$data[19,10] = 0.6094232349989994
# Row 22 (A22=20)
$data[20,0] = 3.636478058440048
$data[20,1] = 0.1864805524337498
$data[20,2] = 0.4867451843504824
$data[20,3] = 0.1346791889329229
$data[20,4] = 0
$data[20,5] = 4.046554680358383
$data[20,6] = 2.850350519552478
$data[20,7] = 0
$data[20,8] = 0.0425056632534595
$data[20,9] = 0
$data[20,10] = 0.6163778447107973
# Row 23 (A23=21)
$data[21,0] = 3.589100630513656
$data[21,1] = 0.1797441614265836
$data[21,2] = 0.4855221218375192
$data[21,3] = 0.1347330474698634
$data[21,4] = 0
$data[21,5] = 4.016620541295083
$data[21,6] = 2.838931294206759
$data[21,7] = 0
$data[21,8] = 0.04241347200254353
$data[21,9] = 0
$data[21,10] = 0.6126410354604417
# Row 24 (A24=22)
$data[22,0] = 3.413314142239244
$data[22,1] = 0.1543061686831493
$data[22,2] = 0.4812406754008123
$data[22,3] = 0.1350013654573257
$data[22,4] = 0
$data[22,5] = 3.906607324823682
$data[22,6] = 2.79777610923486
$data[22,7] = 0
$data[22,8] = 0.04205749706993878
$data[22,9] = 0
$data[22,10] = 0.5990371358311819
# Row 25 (A25=23)
$data[23,0] = 3.231459051957017
$data[23,1] = 0.1270206434398631
$data[23,2] = 0.4773692859357226
$data[23,3] = 0.1354281985194472
$data[23,4] = 0
$data[23,5] = 3.794993238323713
$data[23,6] = 2.757804718985653
$data[23,7] = 0
$data[23,8] = 0.04165878255760624
$data[23,9] = 0
$data[23,10] = 0.5855315558497978

$ws.Range("B2:L25").Value = $data
